# Apply postgame hitter report updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block: rows 10-17 (Inning 3 -> 1) ---
$ws.Range("J10").Value = 1
$ws.Range("M10").Value = ""
$ws.Range("J11").Value = 1
$ws.Range("M12").Value = ""
$ws.Range("M14").Value = "Ground Ball"
$ws.Range("J17").Value = "FB,CB,CH"

# --- Block: rows 19-26 (Inning 1 -> 3) ---
$ws.Range("J19").Value = 3
$ws.Range("M19").Value = ""
$ws.Range("J20").Value = 2
$ws.Range("M21").Value = ""
$ws.Range("M23").Value = "Line Drive"
$ws.Range("J26").Value = "FB,CB,CH"

# --- Block: rows 28-35 (Inning 9 -> 4) ---
$ws.Range("J28").Value = 4
$ws.Range("M28").Value = ""
$ws.Range("J29").Value = 1
$ws.Range("M30").Value = ""
$ws.Range("J32").Value = "Roblez"
$ws.Range("M32").Value = "Fly Ball"
$ws.Range("J33").Value = "Right"
$ws.Range("M33").Value = "Single"
$ws.Range("J34").Value = "88-90 MPH"
$ws.Range("J35").Value = "FB,CB,CH"

# --- Block: rows 37-44 (Inning 8 -> 5) ---
$ws.Range("J37").Value = 5
$ws.Range("M37").Value = ""
$ws.Range("J38").Value = 2
$ws.Range("M39").Value = ""
$ws.Range("J41").Value = "Herbst"
$ws.Range("M41").Value = "Popup"
$ws.Range("J42").Value = "Right"
$ws.Range("M42").Value = "Out"
$ws.Range("J43").Value = "83-85 MPH"
$ws.Range("J44").Value = "SL,FB,CB,CH"

# --- Block: rows 46-48 (Inning 7, only Exit Velo/Launch Angle cleared) ---
$ws.Range("M46").Value = ""
$ws.Range("M48").Value = ""

# --- Block: rows 61-68 (Inning 5 -> 8) ---
$ws.Range("J61").Value = 8
$ws.Range("M61").Value = ""
$ws.Range("J62").Value = 0
$ws.Range("M63").Value = ""
$ws.Range("J65").Value = "Thompson"
$ws.Range("M65").Value = ""
$ws.Range("J66").Value = "Left"
$ws.Range("M66").Value = "Undefined"
$ws.Range("J67").Value = "84-84 MPH"
$ws.Range("J68").Value = "SL,FB,CH"

# --- Block: rows 70-77 (Inning 4 -> 9) ---
$ws.Range("J70").Value = 9
$ws.Range("M70").Value = ""
$ws.Range("J71").Value = 2
$ws.Range("M72").Value = ""
$ws.Range("J74").Value = "Thompson"
$ws.Range("M74").Value = ""
$ws.Range("J75").Value = "Left"
$ws.Range("M75").Value = "Undefined"
$ws.Range("J76").Value = "84-84 MPH"
$ws.Range("J77").Value = "SL,FB,CH"

Write-Host "Applied postgame hitter report updates"
